# [Refactor] 변수명 통일 중 - tier를 tierNum(정수), tierName(문자열), itemClass(서버)로 통일
#
# The underlying grid content (A1:I8) is unchanged except the "tier"
# header in F1, which is renamed to "tier_num". Everything else in the
# source diff is shared-string-table reshuffling that round-trips to the
# same rendered values, plus a couple of view/layout tweaks (active
# selection cell, column F width) that we reproduce here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header F1: "tier" -> "tier_num"
$ws.Range("F1").Value = "tier_num"

# Widen column F (was ~3.8 "chars", now ~8.5 "chars") so the longer
# header text fits; ColumnWidth is expressed in character units.
$ws.Columns.Item(6).ColumnWidth = 7.7

# Move the active selection from D11 to H11, matching the saved
# sheetView selection in the workbook.
[void]$ws.Range("H11").Select()
